# Updated cryptos list (Price and Volume(1h) columns) for rows 2-49,
# and swap/update the NEARProtocol / Decentraland rows (50-51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to store the literal text (avoids Excel auto-converting
    # numeric-looking strings like "1.000" or "16.80" into numbers and
    # losing the exact formatting), while keeping the cell's default style.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2, 4) "28.733.77"
Set-TextCell $ws.Cells.Item(2, 5) "  +6.92%  "
Set-TextCell $ws.Cells.Item(3, 4) "1.808.28"
Set-TextCell $ws.Cells.Item(3, 5) "  +4.41%  "
Set-TextCell $ws.Cells.Item(4, 4) "1.000"
Set-TextCell $ws.Cells.Item(4, 5) "  +0.31%  "
Set-TextCell $ws.Cells.Item(5, 4) "251.67"
Set-TextCell $ws.Cells.Item(5, 5) "  +3.85%  "
Set-TextCell $ws.Cells.Item(6, 4) "1.000"
Set-TextCell $ws.Cells.Item(6, 5) "  +0.27%  "
Set-TextCell $ws.Cells.Item(7, 4) "0.4961"
Set-TextCell $ws.Cells.Item(7, 5) "  +0.92%  "
Set-TextCell $ws.Cells.Item(8, 4) "0.2793"
Set-TextCell $ws.Cells.Item(8, 5) "  +7.29%  "
Set-TextCell $ws.Cells.Item(9, 4) "0.06384"
Set-TextCell $ws.Cells.Item(9, 5) "  +2.52%  "
Set-TextCell $ws.Cells.Item(10, 4) "1.802.41"
Set-TextCell $ws.Cells.Item(10, 5) "  +4.09%  "
Set-TextCell $ws.Cells.Item(11, 4) "16.80"
Set-TextCell $ws.Cells.Item(11, 5) "  +4.32%  "
Set-TextCell $ws.Cells.Item(12, 4) "0.07105"
Set-TextCell $ws.Cells.Item(12, 5) "  +2.86%  "
Set-TextCell $ws.Cells.Item(13, 4) "0.6465"
Set-TextCell $ws.Cells.Item(13, 5) "  +5.63%  "
Set-TextCell $ws.Cells.Item(14, 4) "4.696"
Set-TextCell $ws.Cells.Item(14, 5) "  +4.05%  "
Set-TextCell $ws.Cells.Item(15, 4) "81.99"
Set-TextCell $ws.Cells.Item(15, 5) "  +5.96%  "
Set-TextCell $ws.Cells.Item(16, 4) "28.730.33"
Set-TextCell $ws.Cells.Item(16, 5) "  +6.96%  "
Set-TextCell $ws.Cells.Item(17, 4) "0.9999"
Set-TextCell $ws.Cells.Item(17, 5) "  +0.15%  "
Set-TextCell $ws.Cells.Item(18, 4) "0.000007344"
Set-TextCell $ws.Cells.Item(18, 5) "  +2.18%  "
Set-TextCell $ws.Cells.Item(19, 4) "1.000"
Set-TextCell $ws.Cells.Item(19, 5) "  +0.33%  "
Set-TextCell $ws.Cells.Item(20, 4) "12.26"
Set-TextCell $ws.Cells.Item(20, 5) "  +6.72%  "
Set-TextCell $ws.Cells.Item(21, 4) "2.037.69"
Set-TextCell $ws.Cells.Item(21, 5) "  +4.01%  "
Set-TextCell $ws.Cells.Item(22, 4) "4.606"
Set-TextCell $ws.Cells.Item(22, 5) "  +3.65%  "
Set-TextCell $ws.Cells.Item(23, 4) "8.865"
Set-TextCell $ws.Cells.Item(23, 5) "  +3.32%  "
Set-TextCell $ws.Cells.Item(24, 4) "5.309"
Set-TextCell $ws.Cells.Item(24, 5) "  +3.31%  "
Set-TextCell $ws.Cells.Item(25, 4) "142.70"
Set-TextCell $ws.Cells.Item(25, 5) "  +2.86%  "
Set-TextCell $ws.Cells.Item(26, 4) "16.02"
Set-TextCell $ws.Cells.Item(26, 5) "  +4.50%  "
Set-TextCell $ws.Cells.Item(27, 4) "1.881"
Set-TextCell $ws.Cells.Item(27, 5) "  +5.01%  "
Set-TextCell $ws.Cells.Item(28, 4) "111.40"
Set-TextCell $ws.Cells.Item(28, 5) "  +4.96%  "
Set-TextCell $ws.Cells.Item(29, 4) "1.387"
Set-TextCell $ws.Cells.Item(29, 5) "  +0.54%  "
Set-TextCell $ws.Cells.Item(30, 4) "4.183"
Set-TextCell $ws.Cells.Item(30, 5) "  +5.99%  "
Set-TextCell $ws.Cells.Item(31, 4) "0.08356"
Set-TextCell $ws.Cells.Item(31, 5) "  +4.42%  "
Set-TextCell $ws.Cells.Item(32, 4) "3.836"
Set-TextCell $ws.Cells.Item(32, 5) "  +4.11%  "
Set-TextCell $ws.Cells.Item(33, 4) "0.04942"
Set-TextCell $ws.Cells.Item(33, 5) "  +9.05%  "
Set-TextCell $ws.Cells.Item(34, 4) "1.093"
Set-TextCell $ws.Cells.Item(34, 5) "  +8.42%  "
Set-TextCell $ws.Cells.Item(35, 4) "0.6718"
Set-TextCell $ws.Cells.Item(35, 5) "  +7.41%  "
Set-TextCell $ws.Cells.Item(36, 4) "2.663"
Set-TextCell $ws.Cells.Item(36, 5) "  +2.29%  "
Set-TextCell $ws.Cells.Item(37, 4) "0.9628"
Set-TextCell $ws.Cells.Item(37, 5) "  +2.46%  "
Set-TextCell $ws.Cells.Item(38, 4) "2.638"
Set-TextCell $ws.Cells.Item(38, 5) "  +7.79%  "
Set-TextCell $ws.Cells.Item(39, 4) "2.150"
Set-TextCell $ws.Cells.Item(39, 5) "  +5.00%  "
Set-TextCell $ws.Cells.Item(40, 4) "0.01596"
Set-TextCell $ws.Cells.Item(40, 5) "  +6.05%  "
Set-TextCell $ws.Cells.Item(41, 4) "5.918"
Set-TextCell $ws.Cells.Item(41, 5) "  +4.79%  "
Set-TextCell $ws.Cells.Item(42, 4) "0.9998"
Set-TextCell $ws.Cells.Item(42, 5) "  -0.04%  "
Set-TextCell $ws.Cells.Item(43, 4) "101.27"
Set-TextCell $ws.Cells.Item(43, 5) "  +1.57%  "
Set-TextCell $ws.Cells.Item(44, 4) "0.4115"
Set-TextCell $ws.Cells.Item(44, 5) "  +6.22%  "
Set-TextCell $ws.Cells.Item(45, 4) "7.230"
Set-TextCell $ws.Cells.Item(45, 5) "  +4.22%  "
Set-TextCell $ws.Cells.Item(46, 4) "0.1227"
Set-TextCell $ws.Cells.Item(46, 5) "  +5.58%  "
Set-TextCell $ws.Cells.Item(47, 4) "0.05492"
Set-TextCell $ws.Cells.Item(47, 5) "  +1.94%  "
Set-TextCell $ws.Cells.Item(48, 4) "8.191"
Set-TextCell $ws.Cells.Item(48, 5) "  +2.99%  "
Set-TextCell $ws.Cells.Item(49, 4) "31.33"
Set-TextCell $ws.Cells.Item(49, 5) "  +3.46%  "

# Row 50 becomes NEARProtocol (was Decentraland), row 51 becomes Decentraland (was NEARProtocol)
$ws.Cells.Item(50, 2).Value = "NEARProtocol"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Cells.Item(50, 4) "1.301"
Set-TextCell $ws.Cells.Item(50, 5) "  +4.38%  "

$ws.Cells.Item(51, 2).Value = "Decentraland"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextCell $ws.Cells.Item(51, 4) "0.3605"
Set-TextCell $ws.Cells.Item(51, 5) "  +6.21%  "
